# Append the 06/15/2025 buy to the bitcoin_buys sheet (row 23).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A in this sheet stores the date as literal text (e.g. "06/11/2025"
# in row 22), not as a real date serial. Using Range.Value directly would
# let Excel auto-parse the "06/15/2025" string into a date number, so we
# build it as a text formula first and then paste-special the computed
# value back over itself, which locks it in as a plain text cell without
# touching/creating any NumberFormat-derived cell style.
$ws.Range("A23").Formula = "=""06/15/2025"""
$ws.Range("A23").Copy()
$ws.Range("A23").PasteSpecial(-4163)  # xlPasteValues

$ws.Range("B23").Value = 0.0004702200000000004
$ws.Range("C23").Value = 106333.2057334864
$ws.Range("D23").Value = 50
